$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update reshuffled match rows (id/date/fixed columns unchanged; odds data corrected) ---
# Row 93
$ws.Range("B93").Value = 6236611
$ws.Range("C93").Value = "Venezuela Primera Division"
$ws.Range("D93").Value = "Venezuela Primera Division"
$ws.Range("E93").Value = 45199.6875
$ws.Range("F93").Value = "Mineros"
$ws.Range("G93").Value = "Monagas"
$ws.Range("H93").Value = 2
$ws.Range("I93").Value = 1
$ws.Range("J93").Value = "H"
$ws.Range("K93").Value = 3.2
$ws.Range("L93").Value = 3.4
$ws.Range("M93").Value = 2
$ws.Range("N93").Value = 4.2
$ws.Range("O93").Value = 3.8
$ws.Range("P93").Value = 1.65
$ws.Range("Q93").Value = 0.75
$ws.Range("R93").Value = 1.95
$ws.Range("S93").Value = 1.85
$ws.Range("T93").Value = 2.5
$ws.Range("U93").Value = 1.825
$ws.Range("V93").Value = 1.975
$ws.Range("W93").Value = 3.2
$ws.Range("X93").Value = -1
$ws.Range("Y93").Value = -1
$ws.Range("Z93").Value = 0.95
$ws.Range("AA93").Value = -1
$ws.Range("AB93").Value = 0.825
$ws.Range("AC93").Value = -1

# Row 94
$ws.Range("B94").Value = 6236255
$ws.Range("C94").Value = "Venezuela Primera Division"
$ws.Range("D94").Value = "Venezuela Primera Division"
$ws.Range("E94").Value = 45199.6875
$ws.Range("F94").Value = "Deportivo Rayo Zuliano"
$ws.Range("G94").Value = "Caracas"
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = "D"
$ws.Range("K94").Value = 3.75
$ws.Range("L94").Value = 3.1
$ws.Range("M94").Value = 1.95
$ws.Range("N94").Value = 2.9
$ws.Range("O94").Value = 2.875
$ws.Range("P94").Value = 2.45
$ws.Range("Q94").Value = 0.25
$ws.Range("R94").Value = 1.775
$ws.Range("S94").Value = 2.025
$ws.Range("T94").Value = 2.25
$ws.Range("U94").Value = 1.85
$ws.Range("V94").Value = 1.95
$ws.Range("W94").Value = -1
$ws.Range("X94").Value = 1.875
$ws.Range("Y94").Value = -1
$ws.Range("Z94").Value = 0.3875
$ws.Range("AA94").Value = -0.5
$ws.Range("AB94").Value = -1
$ws.Range("AC94").Value = 0.95

# Row 96
$ws.Range("B96").Value = 6236251
$ws.Range("C96").Value = "Venezuela Primera Division"
$ws.Range("D96").Value = "Venezuela Primera Division"
$ws.Range("E96").Value = 45199.6875
$ws.Range("F96").Value = "Angostura FC"
$ws.Range("G96").Value = "Portuguesa"
$ws.Range("H96").Value = 1
$ws.Range("I96").Value = 2
$ws.Range("J96").Value = "A"
$ws.Range("K96").Value = 3.1
$ws.Range("L96").Value = 3.2
$ws.Range("M96").Value = 2.15
$ws.Range("N96").Value = 4
$ws.Range("O96").Value = 3.6
$ws.Range("P96").Value = 1.75
$ws.Range("Q96").Value = 0.75
$ws.Range("R96").Value = 1.8
$ws.Range("S96").Value = 2
$ws.Range("T96").Value = 2.5
$ws.Range("U96").Value = 1.95
$ws.Range("V96").Value = 1.85
$ws.Range("W96").Value = -1
$ws.Range("X96").Value = -1
$ws.Range("Y96").Value = 0.75
$ws.Range("Z96").Value = -0.5
$ws.Range("AA96").Value = 0.5
$ws.Range("AB96").Value = 0.95
$ws.Range("AC96").Value = -1

# Row 97
$ws.Range("B97").Value = 6236612
$ws.Range("C97").Value = "Venezuela Primera Division"
$ws.Range("D97").Value = "Venezuela Primera Division"
$ws.Range("E97").Value = 45199.6875
$ws.Range("F97").Value = "Zamora"
$ws.Range("G97").Value = "Carabobo"
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 2
$ws.Range("J97").Value = "A"
$ws.Range("K97").Value = 3.2
$ws.Range("L97").Value = 3.1
$ws.Range("M97").Value = 2.15
$ws.Range("N97").Value = 4.5
$ws.Range("O97").Value = 3.3
$ws.Range("P97").Value = 1.75
$ws.Range("Q97").Value = 0.5
$ws.Range("R97").Value = 2
$ws.Range("S97").Value = 1.8
$ws.Range("T97").Value = 2.25
$ws.Range("U97").Value = 1.925
$ws.Range("V97").Value = 1.875
$ws.Range("W97").Value = -1
$ws.Range("X97").Value = -1
$ws.Range("Y97").Value = 0.75
$ws.Range("Z97").Value = -1
$ws.Range("AA97").Value = 0.8
$ws.Range("AB97").Value = -0.5
$ws.Range("AC97").Value = 0.4375

# Row 98
$ws.Range("B98").Value = 6236252
$ws.Range("C98").Value = "Venezuela Primera Division"
$ws.Range("D98").Value = "Venezuela Primera Division"
$ws.Range("E98").Value = 45199.6875
$ws.Range("F98").Value = "Deportivo Tachira"
$ws.Range("G98").Value = "CD Hermanos Colmenares"
$ws.Range("H98").Value = 1
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = "H"
$ws.Range("K98").Value = 1.363
$ws.Range("L98").Value = 4.2
$ws.Range("M98").Value = 7.5
$ws.Range("N98").Value = 1.333
$ws.Range("O98").Value = 4.5
$ws.Range("P98").Value = 8
$ws.Range("Q98").Value = -1.5
$ws.Range("R98").Value = 2
$ws.Range("S98").Value = 1.8
$ws.Range("T98").Value = 2.5
$ws.Range("U98").Value = 1.925
$ws.Range("V98").Value = 1.875
$ws.Range("W98").Value = 0.333
$ws.Range("X98").Value = -1
$ws.Range("Y98").Value = -1
$ws.Range("Z98").Value = -1
$ws.Range("AA98").Value = 0.8
$ws.Range("AB98").Value = -1
$ws.Range("AC98").Value = 0.875

# Row 100
$ws.Range("B100").Value = 6236614
$ws.Range("C100").Value = "Venezuela Primera Division"
$ws.Range("D100").Value = "Venezuela Primera Division"
$ws.Range("E100").Value = 45205.70833333334
$ws.Range("F100").Value = "Mineros"
$ws.Range("G100").Value = "Angostura FC"
$ws.Range("H100").Value = 1
$ws.Range("I100").Value = 2
$ws.Range("J100").Value = "A"
$ws.Range("K100").Value = 2.45
$ws.Range("L100").Value = 3.3
$ws.Range("M100").Value = 2.55
$ws.Range("N100").Value = 1.8
$ws.Range("O100").Value = 3.75
$ws.Range("P100").Value = 3.6
$ws.Range("Q100").Value = -0.5
$ws.Range("R100").Value = 1.825
$ws.Range("S100").Value = 1.975
$ws.Range("T100").Value = 2.75
$ws.Range("U100").Value = 1.8
$ws.Range("V100").Value = 2
$ws.Range("W100").Value = -1
$ws.Range("X100").Value = -1
$ws.Range("Y100").Value = 2.6
$ws.Range("Z100").Value = -1
$ws.Range("AA100").Value = 0.9750000000000001
$ws.Range("AB100").Value = 0.4
$ws.Range("AC100").Value = -0.5

# Row 101
$ws.Range("B101").Value = 6236257
$ws.Range("C101").Value = "Venezuela Primera Division"
$ws.Range("D101").Value = "Venezuela Primera Division"
$ws.Range("E101").Value = 45205.70833333334
$ws.Range("F101").Value = "CD Hermanos Colmenares"
$ws.Range("G101").Value = "Zamora"
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 2
$ws.Range("J101").Value = "A"
$ws.Range("K101").Value = 2.3
$ws.Range("L101").Value = 3.2
$ws.Range("M101").Value = 2.8
$ws.Range("N101").Value = 1.666
$ws.Range("O101").Value = 3.8
$ws.Range("P101").Value = 4.2
$ws.Range("Q101").Value = -0.75
$ws.Range("R101").Value = 1.9
$ws.Range("S101").Value = 1.9
$ws.Range("T101").Value = 2.75
$ws.Range("U101").Value = 1.9
$ws.Range("V101").Value = 1.9
$ws.Range("W101").Value = -1
$ws.Range("X101").Value = -1
$ws.Range("Y101").Value = 3.2
$ws.Range("Z101").Value = -1
$ws.Range("AA101").Value = 0.8999999999999999
$ws.Range("AB101").Value = -1
$ws.Range("AC101").Value = 0.8999999999999999

# --- Append new match row 146 ---
$ws.Range("A146").Value = 144
$ws.Range("B146").Value = 7859982
$ws.Range("C146").Value = "Venezuela Primera Division"
$ws.Range("D146").Value = "Venezuela Primera Division"
$ws.Range("E146").Value = 45347.83333333334
$ws.Range("F146").Value = "Metropolitanos FC"
$ws.Range("G146").Value = "CD Hermanos Colmenares"
$ws.Range("H146").Value = 1
$ws.Range("I146").Value = 2
$ws.Range("J146").Value = "A"
$ws.Range("K146").Value = 1.727
$ws.Range("L146").Value = 3.3
$ws.Range("M146").Value = 4.5
$ws.Range("N146").Value = 1.6
$ws.Range("O146").Value = 3.5
$ws.Range("P146").Value = 5
$ws.Range("Q146").Value = -0.75
$ws.Range("R146").Value = 1.825
$ws.Range("S146").Value = 1.975
$ws.Range("T146").Value = 2.5
$ws.Range("U146").Value = 1.85
$ws.Range("V146").Value = 1.95
$ws.Range("W146").Value = -1
$ws.Range("X146").Value = -1
$ws.Range("Y146").Value = 4
$ws.Range("Z146").Value = -1
$ws.Range("AA146").Value = 0.9750000000000001
$ws.Range("AB146").Value = 0.8500000000000001
$ws.Range("AC146").Value = -1

# Match formatting of preceding data row (bold/bordered id column, date-formatted date column)
$ws.Range("A145:AC145").Copy()
$ws.Range("A146").PasteSpecial(-4122)
$excel.CutCopyMode = 0
